$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Move the "OpenTBS demo" title text from row 2 down to the paragraph
#    below, and put a new, shorter title in its place is NOT what happens;
#    actually B2 keeps the big-title style but now shows "OpenTBS demo"
#    (it used to show the "Far all those raisons..." text because the
#    shared-string table was reordered upstream). We just set the text.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "OpenTBS demo"

# Paragraph block (rows 12-16): re-written text, new bold/colored style.
$ws.Range("B12").Value = "You may consider the following before building your own Microsoft Excel template:"
$ws.Range("B13").Value = "Merging Microsoft Excel templates with OpenTBS has several limitations because of the OpenXML format for Excel."
$ws.Range("B14").Value = "* Formulas won't work because OpenTBS needs to convert cell positions from aboslute to relative in order to have a constistent merged sheet."
$ws.Range("B15").Value = "* Formulas may also make troubles because they are saved twice in the sheet:  one for the expression, and one for the instant result."
$ws.Range("B16").Value = "* Changing picture (using ope=changepic)  because drawing information are saved in another XML sub-file."

# Row 17 used to hold a paragraph; it is now unused -> clear it completely.
$ws.Range("B17").Clear()

# New "Example #1" heading + demo merge-block table (rows 18, 20, 21).
$ws.Range("B18").Value = "Example #1: merging data with rows"

$ws.Range("B20").Value = "First Name"
$ws.Range("C20").Value = "Name"
$ws.Range("D20").Value = "Membership number"

$ws.Range("B21").Value = "[a.firstname;block=row]"
$ws.Range("C21").Value = "[a.name]"
$ws.Range("D21").Value = "[a.number]"

# Old content that used to live at rows 20-21 (single column) is superseded
# by the block above; nothing else to clear since we overwrote B20/B21 and
# C20/D20/C21/D21 were previously empty.

# ---------------------------------------------------------------------------
# Styling
# ---------------------------------------------------------------------------

# B2: big Times New Roman title (unchanged look, just reapplied defensively)
$ws.Range("B2").Font.Name = "Times New Roman"
$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Size = 16
$ws.Range("B2").Font.ColorIndex = 1

# Paragraph block: bold, Accent2-colored text
$paraRange = $ws.Range("B12:B16")
$paraRange.Font.Bold = $true
$paraRange.Font.ThemeColor = 6
$paraRange.Font.TintAndShade = -0.249977111117893

# "Example #1" heading: bold, normal (dark1/theme) text
$ws.Range("B18").Font.Bold = $true
$ws.Range("B18").Font.ThemeColor = 1

# Table header row: bold? no - just shaded fill + thin border all round
$headerRange = $ws.Range("B20:D20")
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.TintAndShade = -0.14999847407452621
$headerRange.Borders.LineStyle = 1

# Table data row: thin border all round, no fill
$dataRange = $ws.Range("B21:D21")
$dataRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths for the new table columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.8
$ws.Columns.Item(3).ColumnWidth = 12.2
$ws.Columns.Item(4).ColumnWidth = 19.0

# ---------------------------------------------------------------------------
# View: selection cursor now rests on B16 instead of B22
# ---------------------------------------------------------------------------
$ws.Range("B16").Select()

Write-Output "edit complete"
